{"js": "// Locate the \"Test Case ID#:  317\" run pair in the table and turn the\n// visible number \"317\" into \"Test_317_07_02\" by inserting \"Test_\" right\n// before it and \"_07_02\" right after it (matching bold formatting is\n// inherited automatically since insertText splits/extends the existing\n// run at the insertion point).\nconst idResults = context.document.body.search(\"317\", { matchCase: true });\nidResults.load(\"text,items\");\nawait context.sync();\n\nif (idResults.items.length === 0) {\n  throw new Error('Could not find \"317\" text to update.');\n}\n\nconst target = idResults.items[0];\ntarget.insertText(\"Test_\", Word.InsertLocation.before);\ntarget.insertText(\"_07_02\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Turn the \"Test Case ID#:  317\" run into \"Test Case ID#:  Test_317_07_02\"\n# by replacing the visible test number \"317\" with \"Test_317_07_02\" using\n# Word's Find/Replace (bold formatting is preserved automatically because\n# Replace only swaps the matched text, leaving the surrounding run\n# properties intact).\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = \"317\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Test_317_07_02\"\n\n$wdFindContinue = 1\n$wdReplaceOne = 1\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceOne) | Out-Null\n"}
